$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("237:238").Insert()

$ws.Range("A237").Value = 11
$ws.Range("B237").Value = "Vega Monumental Concepción"
$ws.Range("C237").Value = "Bíobío"
$ws.Range("D237").Value = 44992
$ws.Range("E237").Value = 8
$ws.Range("F237").Value = 100114013
$ws.Range("G237").Value = "Zanahoria"
$ws.Range("H237").Value = "Sin especificar"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 800
$ws.Range("K237").Value = 5000
$ws.Range("L237").Value = 6000
$ws.Range("M237").Value = 5500
$ws.Range("N237").Value = "`$/saco 20 kilos"
$ws.Range("O237").Value = "Región de Ñuble"
$ws.Range("P237").Value = 275
$ws.Range("Q237").Value = 20
$ws.Range("R237").Value = "Hortaliza"

$ws.Range("A238").Value = 11
$ws.Range("B238").Value = "Vega Monumental Concepción"
$ws.Range("C238").Value = "Bíobío"
$ws.Range("D238").Value = 44992
$ws.Range("E238").Value = 8
$ws.Range("F238").Value = 100114013
$ws.Range("G238").Value = "Zanahoria"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Segunda"
$ws.Range("J238").Value = 400
$ws.Range("K238").Value = 4000
$ws.Range("L238").Value = 4000
$ws.Range("M238").Value = 4000
$ws.Range("N238").Value = "`$/saco 20 kilos"
$ws.Range("O238").Value = "Región de Ñuble"
$ws.Range("P238").Value = 200
$ws.Range("Q238").Value = 20
$ws.Range("R238").Value = "Hortaliza"
